# Generate Report for Handback
# Adds a new handback entry (6712697d-9d87-496c-9440-ea515670d7d5) as row 4
# on the "Overview", "zh-cn" and "de-de" worksheets.

$wb = $excel.ActiveWorkbook

$newGuid = "6712697d-9d87-496c-9440-ea515670d7d5"
$newHash = "ed039b8b219eb22c544bedb23fb79a83e66dec65"

$mdName       = "$newGuid.md"
$zhXlfName    = "$newGuid.$newHash.zh-cn.xlf"
$deXlfName    = "$newGuid.$newHash.de-de.xlf"
$inSyncStatus = "Handed back: in sync with en-US"
$includeReason = "Include"

$zhHandoffDate   = "2016-02-24 06:56:24"
$zhHandbackDate  = "2016-02-24 06:57:16"
$deHandoffDate   = "2016-02-24 06:56:37"
$deHandbackDate  = "2016-02-24 06:57:40"

$hyperlinkColor = 15570276   # matches the workbook's custom hyperlink font (FF6495ED)

function Style-AsHyperlink($range) {
    $range.Font.Underline = 2
    $range.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$mdName", "", "", $mdName)
Style-AsHyperlink $wsOverview.Range("A4")

$wsOverview.Range("B4").Value = $inSyncStatus
$wsOverview.Range("C4").Value = $inSyncStatus

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/master/e2e/$mdName", "", "", $mdName)
Style-AsHyperlink $wsZh.Range("A4")

$wsZh.Range("B4").Value = $inSyncStatus

$wsZh.Hyperlinks.Add($wsZh.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlfName", "", "", $zhXlfName)
Style-AsHyperlink $wsZh.Range("C4")

$wsZh.Range("D4").Value = $zhHandoffDate
$wsZh.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZh.Hyperlinks.Add($wsZh.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/master/e2e/$mdName", "", "", $mdName)
Style-AsHyperlink $wsZh.Range("E4")

$wsZh.Hyperlinks.Add($wsZh.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/master/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlfName", "", "", $zhXlfName)
Style-AsHyperlink $wsZh.Range("F4")

$wsZh.Range("G4").Value = $zhHandbackDate
$wsZh.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZh.Range("H4").Value = $includeReason

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/master/e2e/$mdName", "", "", $mdName)
Style-AsHyperlink $wsDe.Range("A4")

$wsDe.Range("B4").Value = $inSyncStatus

$wsDe.Hyperlinks.Add($wsDe.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlfName", "", "", $deXlfName)
Style-AsHyperlink $wsDe.Range("C4")

$wsDe.Range("D4").Value = $deHandoffDate
$wsDe.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDe.Hyperlinks.Add($wsDe.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/master/e2e/$mdName", "", "", $mdName)
Style-AsHyperlink $wsDe.Range("E4")

$wsDe.Hyperlinks.Add($wsDe.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/master/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlfName", "", "", $deXlfName)
Style-AsHyperlink $wsDe.Range("F4")

$wsDe.Range("G4").Value = $deHandbackDate
$wsDe.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDe.Range("H4").Value = $includeReason
